# feat: add 2022-Q1 data
#
# The workbook originally has two sheets: "2021-Q2" (per-fund holdings for
# that quarter) and "总计" (a running summary of quarter -> fund count ->
# total market value).
#
# This change adds a new quarter of per-fund holdings ("2022-Q1"), reusing
# the "总计" sheet's slot (renaming + rewriting it, mirroring the existing
# "2021-Q2" sheet's layout) and appends a brand-new "总计" summary sheet at
# the end with a new row for 2022-Q1 prepended above the old 2021-Q2 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: the old "总计" sheet becomes the new "2022-Q1" per-fund sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "2022-Q1"

# Keep a handle on the header/index cell styles already present on this
# sheet (s="2" in the saved xf) so the new columns/rows match the look of
# the columns that already used it.
$totalSheet.Range("B1").Copy() | Out-Null
$totalSheet.Range("E1:H1").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("A2").Copy() | Out-Null
$totalSheet.Range("A3").PasteSpecial(-4122) | Out-Null

# Header row.
$totalSheet.Range("B1").Value = "基金代码"
$totalSheet.Range("C1").Value = "基金名称"
$totalSheet.Range("D1").Value = "基金规模"
$totalSheet.Range("E1").Value = "股票总仓位"
$totalSheet.Range("F1").Value = "仓位占比"
$totalSheet.Range("G1").Value = "持有市值(亿元)"
$totalSheet.Range("H1").Value = "仓位排名"

# Helper: write a value that *looks* numeric (leading-zero fund codes,
# decimal-looking percentages, …) as genuine text, the way typing it into
# a Text-formatted cell in Excel would, then drop the cell back to the
# workbook's default ("Normal") style so only the cell's *content* type
# (text) sticks around - not a lingering Text number-format.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - 美元 share class.
$totalSheet.Range("A2").Value = 0
Set-TextValue $totalSheet.Range("B2") "003720"
$totalSheet.Range("C2").Value = "易方达标普生物科技指数（QDII-LOF）美元"
Set-TextValue $totalSheet.Range("D2") "2.11"
Set-TextValue $totalSheet.Range("E2") "94.00"
Set-TextValue $totalSheet.Range("F2") "0.96"
Set-TextValue $totalSheet.Range("G2") "0.0203"
$totalSheet.Range("H2").Value = 9

# Row 3 - 人民币 share class.
$totalSheet.Range("A3").Value = 1
Set-TextValue $totalSheet.Range("B3") "161127"
$totalSheet.Range("C3").Value = "易方达标普生物科技指数（QDII-LOF）人民币"
Set-TextValue $totalSheet.Range("D3") "2.11"
Set-TextValue $totalSheet.Range("E3") "94.00"
Set-TextValue $totalSheet.Range("F3") "0.96"
Set-TextValue $totalSheet.Range("G3") "0.0203"
$totalSheet.Range("H3").Value = 9

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" summary sheet after "2022-Q1", with
# the running per-quarter roll-up (newest quarter on top).
# ---------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $totalSheet)
$newTotal.Name = "总计"

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 0.04

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q2"
$newTotal.Range("C3").Value = 2
$newTotal.Range("D3").Value = 0.04

# Match formatting (bold header / bordered index column) used on the
# "2021-Q2" and "2022-Q1" sheets.
$totalSheet.Range("B1:D1").Copy() | Out-Null
$newTotal.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("A2").Copy() | Out-Null
$newTotal.Range("A2:A3").PasteSpecial(-4122) | Out-Null
